# Apply "ravdess" data-organization updates:
#  - samples_retained!row23 gets the full ravdess row of data (was only A23 filled)
#  - neutral!row4 gets a new "calm" entry (ravdess' neutral-mapped "calm" emotion)
#  - active sheet moves to "neutral", selections updated on both touched sheets

$wb = $excel.ActiveWorkbook

# --- samples_retained: fill in the rest of the "ravdess" row (row 23) ---
$retained = $wb.Worksheets.Item("samples_retained")
$retained.Range("B23").Value = "acted"
$retained.Range("C23").Value = 376
$retained.Range("D23").Value = 1512
$retained.Range("E23").Value = 564
$retained.Range("F23").Value = "English"
$retained.Range("H23").Value = 24
$retained.Range("I23").Value = "neutral, calm, happy, sad, angry, fearful, surprise, and disgust"
$retained.Range("J23").Value = "calm -> 0 here, North American (Canadian) English; contains emotional song samples"

# --- neutral: add a new row for ravdess' "calm" label (mapped to neutral) ---
$neutral = $wb.Worksheets.Item("neutral")
$neutral.Range("A4").Value = "calm"
$neutral.Range("B4").Value = "en"
$neutral.Range("C4").Value = "cal"

# --- view/selection state: neutral tab becomes the active tab/sheet ---
[void]$retained.Range("C24").Select()
[void]$neutral.Activate()
[void]$neutral.Range("A5").Select()
